{"js": "// Remove the explicit <w:tblBorders> override from the two tables that\n// currently force visible single-line borders (table styles \"105\" and\n// \"103\"). The table using style \"104\" already has no visible borders and\n// must stay untouched.\nconst tables = context.document.body.tables;\ntables.load(\"items/style\");\nawait context.sync();\n\nconst targetStyles = [\"105\", \"103\"];\nconst edgeLocations = [\n  Word.BorderLocation.top,\n  Word.BorderLocation.left,\n  Word.BorderLocation.bottom,\n  Word.BorderLocation.right,\n  Word.BorderLocation.insideHorizontal,\n  Word.BorderLocation.insideVertical\n];\n\nfor (const table of tables.items) {\n  if (targetStyles.indexOf(table.style) === -1) {\n    continue;\n  }\n  for (const loc of edgeLocations) {\n    table.getBorder(loc).type = Word.BorderType.none;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the visible single-line borders that were explicitly applied to\n# the two tables that currently show them (the LP/Surat Nomor header table\n# and the Batam/date signature table). The big data table already has no\n# border (LineStyle = wdLineStyleNone/nil) and must be left untouched.\n$d = $word.ActiveDocument\n\n$wdLineStyleNone = [Microsoft.Office.Interop.Word.WdLineStyle]::wdLineStyleNone\n$noBorderTokens = @(\"nil\", \"none\", \"0\")\n\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $table = $d.Tables.Item($i)\n    $topBorder = $table.Borders.Item(1)\n    $currentStyle = [string]$topBorder.LineStyle\n\n    if ($noBorderTokens -notcontains $currentStyle) {\n        for ($b = 1; $b -le 6; $b++) {\n            $table.Borders.Item($b).LineStyle = $wdLineStyleNone\n        }\n    }\n}\n"}
